$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.376.22'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.477.37'
$ws.Range('E3').Value = '  +2.83%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '566.36'
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.60'
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.474.55'
$ws.Range('E9').Value = '  +2.78%  '
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.357'
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.54'
$ws.Range('E14').Value = '  +6.56%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.915.70'
$ws.Range('E15').Value = '  +2.74%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.182.33'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('E17').Value = '  +2.46%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.469.79'
$ws.Range('E18').Value = '  +2.81%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.31'
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '341.57'
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.31'
$ws.Range('E21').Value = '  +1.69%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.77'
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.65'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.172'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.51'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.43'
$ws.Range('E28').Value = '  +4.72%  '
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.77'
$ws.Range('E31').Value = '  +5.93%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0₃0802'
$ws.Range('E32').Value = '  +3.82%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '176.34'
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.53'
$ws.Range('E34').Value = '  +8.06%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '402.33'
$ws.Range('E35').Value = '  +10.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.400'
$ws.Range('E36').Value = '  +1.27%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.88'
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.36'
$ws.Range('E39').Value = '  -3.89%  '
$ws.Range('E40').Value = '  +4.91%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.52'
$ws.Range('E42').Value = '  +3.85%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '150.23'
$ws.Range('E43').Value = '  +4.17%  '
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.72'
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0965'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0519'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.10'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0₆0229'
$ws.Range('E51').Value = '  +5.90%  '
